$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 34 content
$ws.Range("A34").Value = "Encode"
$ws.Range("B34").Value = "Introduce ICONV to convert encoding and XXD to cat in HEX"
$ws.Range("C34").Value = "Convert from file a to file b: `n> iconv -f ASCII -t UTF-8 < fileA.txt > fileB.txt`nShow the encoding option:`n> iconv -l `nShow HEX of a file:`n> xxd filea.txt"

# Match formatting of the row above (wrap text style) and row height
$ws.Range("C34").WrapText = $true
$ws.Rows.Item(34).RowHeight = 94.5

# Update view to scroll and select like the author's saved state
$ws.Application.ActiveWindow.ScrollRow = 23
$ws.Range("C35").Select()
